$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - I1 ("I0") and J1 ("IF"), matching the style used by H1 (bordered/bold header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows (I and J columns), unstyled like the rest of the numeric columns
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 4

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 3

$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 7

$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 4
